# Auto-generated: updates cryptos Price (D) / Volume(1h) (E) cells
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.914.66'
$ws.Range('E2').Value = '  +3.85%  '
$ws.Range('D3').Value = '3.712.14'
$ws.Range('E3').Value = '  +8.33%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''590.10'
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('D6').Value = '''180.61'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').Value = '3.703.36'
$ws.Range('E7').Value = '  +8.28%  '
$ws.Range('D8').Value = '''0.618'
$ws.Range('E8').Value = '  +4.27%  '
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('E11').Value = '  +4.92%  '
$ws.Range('D12').Value = '''50.06'
$ws.Range('E12').Value = '  +3.73%  '
$ws.Range('D13').Value = '''0.0000288'
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('D14').Value = '4.312.76'
$ws.Range('E14').Value = '  +8.47%  '
$ws.Range('D15').Value = '''683.32'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('E16').Value = '  +4.40%  '
$ws.Range('D17').Value = '3.714.67'
$ws.Range('E17').Value = '  +8.27%  '
$ws.Range('D18').Value = '72.004.96'
$ws.Range('E18').Value = '  +3.89%  '
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('D21').Value = '''11.71'
$ws.Range('E21').Value = '  +3.26%  '
$ws.Range('D22').Value = '''6.47'
$ws.Range('E22').Value = '  +20.22%  '
$ws.Range('D23').Value = '''0.946'
$ws.Range('E23').Value = '  +3.62%  '
$ws.Range('D24').Value = '''17.89'
$ws.Range('E24').Value = '  +5.43%  '
$ws.Range('D25').Value = '''103.97'
$ws.Range('E25').Value = '  +3.03%  '
$ws.Range('D26').Value = '''4.06'
$ws.Range('E26').Value = '  +4.10%  '
$ws.Range('E27').Value = '  +5.27%  '
$ws.Range('D28').Value = '''10.32'
$ws.Range('E28').Value = '  +6.84%  '
$ws.Range('E29').Value = '  +5.97%  '
$ws.Range('D30').Value = '''9.30'
$ws.Range('E30').Value = '  +6.15%  '
$ws.Range('D31').Value = '''7.39'
$ws.Range('E31').Value = '  +7.35%  '
$ws.Range('E32').Value = '  +14.26%  '
$ws.Range('D33').Value = '''11.31'
$ws.Range('E33').Value = '  +2.44%  '
$ws.Range('E34').Value = '  +4.27%  '
$ws.Range('D35').Value = '''562.99'
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('D36').Value = '''59.74'
$ws.Range('E36').Value = '  +3.16%  '
$ws.Range('D37').Value = '3.770.27'
$ws.Range('E37').Value = '  +3.78%  '
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('E39').Value = '  +2.91%  '
$ws.Range('D40').Value = '0.0₃0782'
$ws.Range('E40').Value = '  +6.55%  '
$ws.Range('D41').Value = '''35.81'
$ws.Range('E41').Value = '  +1.98%  '
$ws.Range('D42').Value = '''3.48'
$ws.Range('E42').Value = '  +6.29%  '
$ws.Range('D43').Value = '''0.0465'
$ws.Range('E43').Value = '  +9.41%  '
$ws.Range('E44').Value = '  +4.25%  '
$ws.Range('E45').Value = '  +5.52%  '
$ws.Range('D46').Value = '''2.91'
$ws.Range('E46').Value = '  +9.07%  '
$ws.Range('D47').Value = '''3.37'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('E48').Value = '  +3.86%  '
$ws.Range('E49').Value = '  +2.23%  '
$ws.Range('D50').Value = '''0.999'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').Value = '''135.51'
$ws.Range('E51').Value = '  +3.50%  '
